$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50; existing rows 50-193 shift down to 51-194
$ws.Rows.Item(50).Insert()

# Fill the new row 50 with the weekly data point (matching the fixed columns used
# throughout this sheet for this market/category, plus the new date/price values)
$ws.Range("A50").Value = 4
$ws.Range("B50").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C50").Value = "Los Lagos"
$ws.Range("D50").Value = 44525
$ws.Range("E50").Value = 10
$ws.Range("F50").Value = 100112037
$ws.Range("G50").Value = "Cebollín"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 70
$ws.Range("K50").Value = 6000
$ws.Range("L50").Value = 6000
$ws.Range("M50").Value = 6000
$ws.Range("N50").Value = "$/paquete 36 unidades"
$ws.Range("O50").Value = "Región Metropolitana"
$ws.Range("P50").Value = 167
$ws.Range("Q50").Value = 36
$ws.Range("R50").Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Range("D50").NumberFormat = $ws.Range("D51").NumberFormat
